$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Health Care Industry KPIs")

# --- Reference link header & url (columns H1:H2) ---
$ws.Range("H1").Value = "Reference Link "
$ws.Range("H2").Value = "https://insightsoftware.com/blog/25-best-healthcare-kpis-and-metric-examples/"

# --- New KPI rows (12-16) ---
$ws.Range("A12").Value = "Staff to Patient Ratio"
$ws.Range("B12").Value = "Ratio b/w staff and patients"

$ws.Range("A13").Value = "Child Immunization Ratio"
$ws.Range("B13").Value = "Child vaccination related KPI"

$ws.Range("A14").Value = "Net Profit Margin"
$ws.Range("B14").Value = "Net profit margin per year"

$ws.Range("A15").Value = "Operating Cashflow"
$ws.Range("B15").Value = " Current Cashflow Rate."

$ws.Range("A16").Value = "Readmission Rate"
$ws.Range("B16").Value = "Readmission of patient to hospital."

# Match the formatting (style) used by the rest of column A for the new rows,
# and give the new "Child Immunization Ratio" row the taller (30pt) row height
# used elsewhere in the sheet for wrapped / longer rows.
$ws.Range("A2").Copy()
$ws.Range("A12:A17").PasteSpecial(-4122)

$ws.Rows.Item(13).RowHeight = 30

# --- Trailing empty, but styled, row 17 (A only) ---
$ws.Range("A17").Value = ""

# --- Page setup (portrait orientation) ---
$ws.PageSetup.Orientation = 1

# --- Selection moves to O16 ---
$null = $ws.Range("O16").Select()

Write-Output "done"
